# Refresh cryptocurrency Price (D) and Volume(1h) (E) figures to match the
# latest coinranking.com snapshot pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.210.20"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "3.918.47"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'486.86"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'146.48"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.733"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").Value = "'0.0000349"
$ws.Range("E11").Value = "  -4.90%  "
$ws.Range("D12").Value = "'43.25"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "'10.71"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "4.544.71"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "3.907.88"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "'14.25"
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "'20.09"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "68.333.86"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'431.35"
$ws.Range("E21").Value = "  -1.41%  "
$ws.Range("E22").Value = "  +3.43%  "
$ws.Range("D23").Value = "'15.16"
$ws.Range("E23").Value = "  +4.46%  "
$ws.Range("D24").Value = "'88.38"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'11.72"
$ws.Range("E25").Value = "  +20.36%  "
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("D27").Value = "'11.19"
$ws.Range("E27").Value = "  +9.42%  "
$ws.Range("D28").Value = "'37.83"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("D30").Value = "'718.20"
$ws.Range("E30").Value = "  -2.14%  "
$ws.Range("D31").Value = "'13.73"
$ws.Range("E31").Value = "  +2.98%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").Value = "'6.19"
$ws.Range("E35").Value = "  +14.76%  "
$ws.Range("D36").Value = "'41.63"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").Value = "'60.98"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "'0.397"
$ws.Range("E38").Value = "  +17.69%  "
$ws.Range("E39").Value = "  -3.68%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "'3.00"
$ws.Range("E41").Value = "  +16.73%  "
$ws.Range("D42").Value = "'0.0492"
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("E44").Value = "  +4.54%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'3.37"
$ws.Range("E46").Value = "  +3.40%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("D50").Value = "'144.84"
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("D51").Value = "0.0₆0339"
$ws.Range("E51").Value = "  +26.33%  "
